# Update "想去人数" (want-to-go count) figures in the 展览 and 全部类型 sheets.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 5360
$ws1.Range("F7").Value = 607
$ws1.Range("F11").Value = 1472
$ws1.Range("F12").Value = 4293
$ws1.Range("F17").Value = 3429

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 5360
$ws4.Range("F8").Value = 607
$ws4.Range("F12").Value = 1472
$ws4.Range("F13").Value = 4293
$ws4.Range("F18").Value = 3429
